$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 7 (ZOTAC GAMING GeForce RTX 3070 / B08LF1CWT2 / 700),
# pushing the existing rows 7-10 down to rows 10-13.
$ws.Range("A7:C9").EntireRow.Insert()

# New row 7: EVGA RTX 3080 FTW3 ULTRA
$ws.Range("A7").Value = "EVGA 10G-P5-3897-KR GeForce RTX 3080 FTW3 ULTRA"
$ws.Range("B7").Value = "B08HR3Y5GQ"
$ws.Range("C7").Value = 830

# New row 8: ASUS ROG STRIX RTX 3080
$ws.Range("A8").Value = "ASUS ROG STRIX NVIDIA GeForce RTX 3080"
$ws.Range("B8").Value = "B08J6F174Z"
$ws.Range("C8").Value = 1150

# New row 9: ASUS TUF Gaming RTX 3080
$ws.Range("A9").Value = "ASUS TUF Gaming NVIDIA GeForce RTX 3080"
$ws.Range("B9").Value = "B08HH5WF97"
$ws.Range("C9").Value = 750

$ws.Range("A14").Select()
